# Applies the "added metro to osc reading" commit to ModMatrixTable.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Main")

# --- Column widths (target raw widths aren't reachable exactly through
# ColumnWidth's MDW pixel-grid rounding here, so use the nearest inputs
# that land in the same rounded bucket as the target) ----------------------
$ws.Columns.Item(2).ColumnWidth = 42.5
$ws.Columns.Item(5).ColumnWidth = 37.833
$ws.Columns.Item(6).ColumnWidth = 14.5
$ws.Columns.Item(7).ColumnWidth = 30.667
$ws.Columns.Item(8).ColumnWidth = 22.0
$ws.Columns.Item(9).ColumnWidth = 22.667
$ws.Columns.Item(10).ColumnWidth = 35.667

# --- Header row (row 2): rename mod-matrix parameter columns -------------
$ws.Range("B2").Value = "/track/1/fx/1/fxparam/1/value"
$ws.Range("C2").Value = "/track/1/fx/1/fxparam/2/value"
$ws.Range("D2").Value = "/track/1/fx/1/fxparam/3/value"
$ws.Range("E2").Value = "/track/1/fx/1/fxparam/4/value"
$ws.Range("F2").Value = "/track/1/fx/2/fxparam/1/value"
$ws.Range("G2").Value = "/track/1/fx/2/fxparam/2/value"
$ws.Range("H2").Value = "/track/1/fx/3/fxparam/3/value"
$ws.Range("I2").Value = "/track/1/fx/4/fxparam/4/value"

# column J (LfoAmp) is no longer used - clear header and match style of C2
$ws.Range("J2").ClearContents()
$ws.Range("B2").Copy()
$ws.Range("C2").PasteSpecial(-4122)
$ws.Range("F2").PasteSpecial(-4122)
$ws.Range("G2").PasteSpecial(-4122)
$ws.Range("D2").Copy()
$ws.Range("H2").PasteSpecial(-4122)
$ws.Range("I2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Clear the now-unused column J data (rows 3-11) -----------------------
$ws.Range("J3:J11").ClearContents()

# --- Selection matches the author's final click (column J) ----------------
$ws.Range("J1:J1048576").Select()
